# bug fixes on test 5 generator
# Updates the "Write Latency" min (O) and max (P) columns for rows 3-23
# on the active worksheet to reflect corrected generator output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;  O = "1378"; P = "2966.9k" },
    @{ Row = 4;  O = "1398"; P = "3559.6k" },
    @{ Row = 5;  O = "1415"; P = "11024k" },
    @{ Row = 6;  O = "1368"; P = "565907" },
    @{ Row = 7;  O = "1421"; P = "4460.0k" },
    @{ Row = 8;  O = "1267"; P = "2499.7k" },
    @{ Row = 9;  O = "1370"; P = "2111.7k" },
    @{ Row = 10; O = "1326"; P = "7710.7k" },
    @{ Row = 11; O = "1326"; P = "2150.2k" },
    @{ Row = 12; O = "2";    P = "3003" },
    @{ Row = 13; O = "1445"; P = "3773.0k" },
    @{ Row = 14; O = "2";    P = "2248" },
    @{ Row = 15; O = "2";    P = "9009" },
    @{ Row = 16; O = "10";   P = "6647" },
    @{ Row = 17; O = "1264"; P = "3352.6k" },
    @{ Row = 18; O = "1438"; P = "1010.8k" },
    @{ Row = 19; O = "1557"; P = "2135.2k" },
    @{ Row = 20; O = "2";    P = "1696" },
    @{ Row = 21; O = "2";    P = "6448" },
    @{ Row = 22; O = "1290"; P = "991327" },
    @{ Row = 23; O = "1409"; P = "2127.7k" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 15).Value = $u.O
    $ws.Cells.Item($u.Row, 16).Value = $u.P
}
